$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Euramet")

# Row 29
$ws.Range("F29").Value = -157.0411782090434

# Row 30
$ws.Range("F30").Value = -156.5579745837849

# Row 31
$ws.Range("F31").Value = -156.5579745837849
$ws.Range("G31").Value = -1.1772

# Row 33
$ws.Range("F33").Value = -156.5579745837849

# Row 34
$ws.Range("F34").Value = -156.5579745837849

# Row 36
$ws.Range("G36").Value = -1.1772

# Row 37
$ws.Range("F37").Value = -156.5579745837849

# Row 38
$ws.Range("G38").Value = -1.1772

# Row 39
$ws.Range("F39").Value = -157.0411782090434

# Row 40
$ws.Range("F40").Value = -156.5579745837849
$ws.Range("G40").Value = -1.3734

# Row 44
$ws.Range("F44").Value = -156.5579745837849
$ws.Range("G44").Value = -1.1772

# Row 45 (previously empty, now populated)
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = -156.5579745837849
$ws.Range("G45").Value = -1.1772
$ws.Range("H45").Value = 1

# Row 46 (previously empty, now populated)
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = -156.5579745837849
$ws.Range("G46").Value = -1.1772
$ws.Range("H46").Value = 1
